$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.338.86'
$ws.Range('E2').Value = '  +4.09%  '
$ws.Range('D3').Value = '2.431.64'
$ws.Range('E3').Value = '  +3.11%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.88'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.44'
$ws.Range('E6').Value = '  +3.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +2.08%  '
$ws.Range('E9').Value = '  +5.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.76'
$ws.Range('E10').Value = '  +3.98%  '
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('E13').Value = '  +5.42%  '
$ws.Range('D14').Value = '2.863.94'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').Value = '60.246.33'
$ws.Range('E15').Value = '  +4.05%  '
$ws.Range('E16').Value = '  +4.15%  '
$ws.Range('D17').Value = '2.447.82'
$ws.Range('E17').Value = '  +3.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.39'
$ws.Range('E18').Value = '  +6.15%  '
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '334.22'
$ws.Range('E20').Value = '  +1.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.77'
$ws.Range('E21').Value = '  +0.76%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.27'
$ws.Range('E23').Value = '  +4.39%  '
$ws.Range('E24').Value = '  +3.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.65'
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.36'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('D28').Value = '0.0₃0791'
$ws.Range('E28').Value = '  +7.17%  '
$ws.Range('E29').Value = '  +1.61%  '
$ws.Range('E30').Value = '  +3.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.37'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.06'
$ws.Range('E32').Value = '  +2.56%  '
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +6.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.24'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '324.18'
$ws.Range('E39').Value = '  +11.68%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '39.82'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.421'
$ws.Range('E41').Value = '  +10.93%  '
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '140.85'
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('E44').Value = '  +3.58%  '
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.68'
$ws.Range('E46').Value = '  +2.67%  '
$ws.Range('E47').Value = '  +9.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.573'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('E51').Value = '  -0.23%  '
